$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.751.42"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "1.803.48"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").Value = "'231.49"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").Value = "'0.5951"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'0.2781"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").Value = "'0.06840"
$ws.Range("D10").Value = "'23.39"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").Value = "'0.07546"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").Value = "1.805.20"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "'4.766"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").Value = "'0.6262"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "2.047.65"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "'0.000009288"
$ws.Range("E16").Value = "  -7.77%  "
$ws.Range("D17").Value = "'75.37"
$ws.Range("E17").Value = "  -4.66%  "
$ws.Range("D18").Value = "28.698.68"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("E19").Value = "  -6.87%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "'210.93"
$ws.Range("E21").Value = "  -7.35%  "
$ws.Range("D22").Value = "'11.45"
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("D23").Value = "'6.863"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").Value = "'154.34"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "'7.852"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("E27").Value = "  -2.62%  "
$ws.Range("D28").Value = "'16.40"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("D30").Value = "'0.06210"
$ws.Range("E30").Value = "  -3.82%  "
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").Value = "'3.776"
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("D33").Value = "'3.756"
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("D35").Value = "'1.062"
$ws.Range("E35").Value = "  -5.68%  "
$ws.Range("D36").Value = "'0.6404"
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("D37").Value = "'2.494"
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("D38").Value = "'2.723"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").Value = "'0.01710"
$ws.Range("E39").Value = "  -1.92%  "
$ws.Range("D40").Value = "'6.426"
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("D41").Value = "1.140.17"
$ws.Range("E41").Value = "  -6.55%  "
$ws.Range("D42").Value = "'0.8675"
$ws.Range("E42").Value = "  -6.62%  "
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").Value = "'100.56"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "1.964.93"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").Value = "'60.61"
$ws.Range("E46").Value = "  -3.83%  "
$ws.Range("E47").Value = "  -5.16%  "
$ws.Range("D48").Value = "'1.599"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "'8.366"
$ws.Range("E49").Value = "  -2.58%  "
$ws.Range("D50").Value = "'0.05472"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").Value = "'0.4491"
$ws.Range("E51").Value = "  -1.51%  "
